$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '96.178.65'
$ws.Range('E2').Value = '  +4.71%  '
$ws.Range('D3').Value = '3.663.82'
$ws.Range('E3').Value = '  +10.12%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '240.68'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +4.53%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '644.23'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.86%  '
$ws.Range('E7').Value = '  +6.29%  '
$ws.Range('E8').Value = '  +5.44%  '
$ws.Range('E9').Value = '  -0.13%  '
$ws.Range('E10').Value = '  +5.75%  '
$ws.Range('D11').Value = '3.661.94'
$ws.Range('E11').Value = '  +9.99%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '43.79'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +2.53%  '
$ws.Range('E13').Value = '  +3.92%  '
$ws.Range('E14').Value = '  +3.53%  '
$ws.Range('D15').Value = '4.349.03'
$ws.Range('E15').Value = '  +10.13%  '
$ws.Range('D16').Value = '95.977.46'
$ws.Range('E16').Value = '  +4.59%  '
$ws.Range('E17').Value = '  +5.88%  '
$ws.Range('E18').Value = '  +4.58%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.61'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +25.91%  '
$ws.Range('D20').Value = '3.657.99'
$ws.Range('E20').Value = '  +9.71%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '18.54'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +6.73%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '518.11'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +5.17%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.46'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.65%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.486'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +10.99%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.0000199'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +9.25%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '6.77'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +5.18%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '97.88'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +8.81%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '12.57'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +5.77%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '3.19'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +22.73%  '
$ws.Range('E30').Value = '  +5.12%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.142'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +3.13%  '
$ws.Range('E32').Value = '  -0.09%  '
$ws.Range('E33').Value = '  +4.88%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.997'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.39%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '32.22'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +13.90%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.580'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +10.34%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '566.19'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.17%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '7.87'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +6.95%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.50'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +9.91%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.950'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +9.63%  '
$ws.Range('E41').Value = '  +2.97%  '
$ws.Range('E42').Value = '  -0.04%  '
$ws.Range('E43').Value = '  +4.30%  '
$ws.Range('E44').Value = '  +4.51%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '5.75'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +6.54%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '23.76'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.33%  '
$ws.Range('E47').Value = '  +6.11%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '54.15'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +4.65%  '
$ws.Range('E49').Value = '  -2.97%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.30'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +4.16%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '208.84'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +16.08%  '
